{"js": "// The document uses curly-quoted placeholders such as\n//   \u201c{{SEXO_12}} PROMITENTE \u2026   and   \u201c{{SEXO_13}} PROMITENTE(S) \u2026\n// The edit removes the opening curly quote (\u201c) that precedes the\n// {{SEXO_12}} / {{SEXO_13}} placeholder and also removes the space\n// between the closing \"}}\" and the following word \"PROMITENTE\",\n// e.g. \u201c{{SEXO_12}} PROMITENTE  ->  {{SEXO_12}}PROMITENTE\n// This occurs at every location in the document where that exact\n// sequence appears (11 occurrences total), so a global search/replace\n// covers every hunk in the diff.\n\nconst body = context.document.body;\n\nconst patterns = [\n  \"\\u201c{{SEXO_12}} PROMITENTE\",\n  \"\\u201c{{SEXO_13}} PROMITENTE\",\n];\n\nfor (const pattern of patterns) {\n  const results = body.search(pattern, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  const replacement = pattern.replace(\"\\u201c\", \"\").replace(\"}} PROMITENTE\", \"}}PROMITENTE\");\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document uses curly-quoted placeholders such as\n#   \"{{SEXO_12}} PROMITENTE ...   and   \"{{SEXO_13}} PROMITENTE(S) ...\n# (opening curly quote U+201C before the placeholder).\n# This edit removes that opening curly quote and also removes the\n# space between the closing \"}}\" and the following word \"PROMITENTE\",\n# e.g.  \"{{SEXO_12}} PROMITENTE  ->  {{SEXO_12}}PROMITENTE\n# The exact sequence occurs 11 times throughout the document (for both\n# SEXO_12 and SEXO_13), so a pair of document-wide Find/Replace passes\n# covers every location touched by the diff.\n\n$d = $word.ActiveDocument\n$openQuote = [char]0x201C\n\n$patterns = @(\"12\", \"13\")\n\nforeach ($num in $patterns) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $openQuote + \"{{SEXO_\" + $num + \"}} PROMITENTE\"\n    $rng.Find.Replacement.Text = \"{{SEXO_\" + $num + \"}}PROMITENTE\"\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1  # wdFindContinue\n    $rng.Find.Execute($rng.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng.Find.Replacement.Text, 2) | Out-Null\n}\n"}
